$d = $word.ActiveDocument

# "REST api and playlist": the videolist field used to carry the raw
# YouTube watch-id together with the playlist query string
# (coIKdmZb6Tw&list=PLAYLIST_ID). Now that the player is driven purely
# off the playlist id (resolved through the REST api), only the
# playlist id itself is kept in the [...] value.
$d.Content.Find.Execute(
    "coIKdmZb6Tw&list=PLE1P7Fwrqlm-l4HO-w26R3-EUuDI0AVAB",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "PLE1P7Fwrqlm-l4HO-w26R3-EUuDI0AVAB",
    2
) | Out-Null

# Word stamps a hidden "_GoBack" bookmark at the position of the last
# edit; drop it so it doesn't linger on the now-edited paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
